$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in hours for Week 7 timesheet
$ws.Range("H6").Value = 0.5
$ws.Range("F8").Value = 1
$ws.Range("G12").Value = 1
$ws.Range("H12").Value = 1

# Recalculate so the SUM() formulas (daily/weekly totals) pick up new values
$excel.CalculateFull()

# Update the active selection to H17, matching where the user left off editing
$ws.Range("H17").Select()
